$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nk")
$ws.Range("A9").Value = "Provola"
